$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Segment names currently sitting in column A (rows 2-20), in row order.
$segments = @(
    "background",
    "back_bumper",
    "back_glass",
    "back_left_door",
    "back_left_light",
    "back_right_door",
    "back_right_light",
    "front_bumper",
    "front_glass",
    "front_left_door",
    "front_left_light",
    "front_right_door",
    "front_right_light",
    "hood",
    "left_mirror",
    "right_mirror",
    "tailgate",
    "trunk",
    "wheel"
)

# Insert a brand new column B - shifts old B..K (PercActivations..totalStd)
# one column to the right, into C..L, and leaves a blank column B behind.
$ws.Columns("B:B").Insert()

# New header cell for the inserted column.
$ws.Range("B1").Value = "segments"
# Give it the same header formatting (bold / centered / bordered) as the
# rest of row 1 by copying the format from the neighboring header cell.
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122)

# Move the segment names out of column A and into the new column B, and
# replace column A with the plain numeric index (0-based).
for ($i = 0; $i -lt $segments.Count; $i++) {
    $row = $i + 2
    $ws.Range("B$row").Value = $segments[$i]
    # The moved-in label should carry no special styling (matches the
    # plain data cells), so strip whatever the Insert step left behind.
    $ws.Range("B$row").ClearFormats()
    $ws.Range("A$row").Value = $i
}

$excel.CutCopyMode = $false
